# cambio para la selección por orden de opciones de usaurio
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fecha (date serial) for all data rows
$newFecha = 45689

# Row data: Col D (cierre_tpv_desc), Col E (Nombre_MdP), Col F (total_ventas), Col G (total_operaciones)
$rows = @(
    @{ Row = 2;  D = "Mañana"; E = "EUROS";        F = 4010.85;  G = 450 },
    @{ Row = 3;  D = "Mañana"; E = "TARJETA VISA"; F = 7933.5;   G = 711 },
    @{ Row = 4;  D = "Tarde";  E = "EUROS";        F = 3100.5;   G = 252 },
    @{ Row = 5;  D = "Tarde";  E = "TARJETA VISA"; F = 7644.6;   G = 630 },
    @{ Row = 6;  D = "Mañana"; E = "EUROS";        F = 7256.88;  G = 693 },
    @{ Row = 7;  D = "Mañana"; E = "TARJETA VISA"; F = 20079.27; G = 1386 },
    @{ Row = 8;  D = "Tarde";  E = "EUROS";        F = 6324.84;  G = 711 },
    @{ Row = 9;  D = "Tarde";  E = "GLOVO";        F = 327.6;    G = 18 },
    @{ Row = 10; D = "Tarde";  E = "SMS";          F = 42.3;     G = 9 },
    @{ Row = 11; D = "Tarde";  E = "TARJETA VISA"; F = 16288.2;  G = 1386 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 3).Value = $newFecha
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
